$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column is stored as text (not auto-converted to a number),
# matching the source data which keeps prices/volumes as formatted strings.
$ws.Range("D2:D15").NumberFormat = "@"
$ws.Range("D17:D33").NumberFormat = "@"
$ws.Range("D34:D51").NumberFormat = "@"

# Update Price (D) and Volume(1h) (E) columns for rows 2-33 (Coin/Link unchanged)
$ws.Range("D2").Value = "26.506.00"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.730.72"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "247.19"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "0.4845"
$ws.Range("E7").Value = "  +0.45%  "
$ws.Range("D8").Value = "0.2669"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").Value = "0.06221"
$ws.Range("E9").Value = "  -0.75%  "
$ws.Range("D10").Value = "1.729.23"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "0.07063"
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").Value = "15.65"
$ws.Range("E12").Value = "  -1.44%  "
$ws.Range("D13").Value = "4.617"
$ws.Range("E13").Value = "  +1.98%  "
$ws.Range("D14").Value = "0.6118"
$ws.Range("E14").Value = "  -2.18%  "
$ws.Range("D15").Value = "77.39"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D17").Value = "26.497.23"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "0.9997"
$ws.Range("E18").Value = "  -0.03%  "
$ws.Range("D19").Value = "0.000007230"
$ws.Range("E19").Value = "  +4.79%  "
$ws.Range("D20").Value = "11.56"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").Value = "1.952.81"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("D22").Value = "4.511"
$ws.Range("E22").Value = "  -2.75%  "
$ws.Range("D23").Value = "8.791"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "5.258"
$ws.Range("E24").Value = "  -2.32%  "
$ws.Range("D25").Value = "137.46"
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "15.44"
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "1.779"
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").Value = "108.25"
$ws.Range("E28").Value = "  +1.34%  "
$ws.Range("D29").Value = "1.397"
$ws.Range("E29").Value = "  -2.53%  "
$ws.Range("D30").Value = "3.985"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "0.07990"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("D32").Value = "3.689"
$ws.Range("E32").Value = "  -1.84%  "
$ws.Range("D33").Value = "0.04564"
$ws.Range("E33").Value = "  -1.57%  "

# Rows 34-51: a new "Frax" row is inserted, shifting the remaining rows down by one
# (the previous last row, Aave, is pushed off the bottom of the table)
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").Value = "0.9996"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.612"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.006"
$ws.Range("E36").Value = "  +0.59%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "0.6344"
$ws.Range("E37").Value = "  -1.59%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "0.8941"
$ws.Range("E38").Value = "  -4.86%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "2.015"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "2.388"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "1.003"
$ws.Range("E41").Value = "  +0.35%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "0.01504"
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "101.46"
$ws.Range("E43").Value = "  -10.32%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.477"
$ws.Range("E44").Value = "  -4.75%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.3897"
$ws.Range("E45").Value = "  -0.71%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").Value = "7.021"
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.1183"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.05382"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.917"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "30.58"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.251"
$ws.Range("E51").Value = "  -1.36%  "
